# Definition of Rounds 1 - 4
# Applies the "Runde 4" / Finale category+winner additions to the
# "Kategorien" sheet, marks the category-group headers in column A with
# the built-in "Explanatory Text" style, and updates the sheet selections.

$wb = $excel.ActiveWorkbook
$wsSpieler = $wb.Worksheets.Item("Spieler")
$wsKategorien = $wb.Worksheets.Item("Kategorien")

# ---------------------------------------------------------------------
# Kategorien sheet: mark the "group header" rows in column A with the
# built-in Explanatory-Text cell style (italic grey) - this is the style
# normally used to tag the question-category group a round's rows belong
# to.
# ---------------------------------------------------------------------
$explanatoryRows = @(2,3,4,5,6,7,8,9,10,13,16,17,19,20)
foreach ($r in $explanatoryRows) {
    $wsKategorien.Range("A$r").Style = "Explanatory Text"
}

# ---------------------------------------------------------------------
# Kategorien sheet: fill in winners (column K) for rounds 1-3 that were
# already defined.
# ---------------------------------------------------------------------
$wsKategorien.Range("K2").Value = "Michael"
$wsKategorien.Range("K3").Value = "Marco"
$wsKategorien.Range("K4").Value = "Christian"
$wsKategorien.Range("K5").Value = "Arkadi"

$wsKategorien.Range("K10").Value = "Miram"
$wsKategorien.Range("K11").Value = "Anton"
$wsKategorien.Range("K12").Value = "Roger"
$wsKategorien.Range("K13").Value = "Adrian"

# ---------------------------------------------------------------------
# Kategorien sheet: Round 3 (row 17 header) - row 14 category got
# corrected from "Recursive Acr"/Technisch to "TV Themes"/Popkultur, and
# rows 18-20 gain their Round-3 category/type/winner entries.
# ---------------------------------------------------------------------
$wsKategorien.Range("I14").Value = "TV Themes"
$wsKategorien.Range("J14").Value = "Popkultur"
$wsKategorien.Range("J14").Style = "Gut"

$wsKategorien.Range("I18").Value = "TLDs"
$wsKategorien.Range("J18").Value = "Technisch"
$wsKategorien.Range("J18").Style = "Schlecht"
$wsKategorien.Range("K18").Value = "Jochen"

$wsKategorien.Range("I19").Value = "Serien Rollen"
$wsKategorien.Range("J19").Value = "Popkultur"
$wsKategorien.Range("J19").Style = "Gut"
$wsKategorien.Range("K19").Value = "Roland"

$wsKategorien.Range("I20").Value = "Virtual"
$wsKategorien.Range("J20").Value = "Popkultur"
$wsKategorien.Range("J20").Style = "Gut"
$wsKategorien.Range("K20").Value = "Björn"

$wsKategorien.Range("I21").Value = "Seven"
$wsKategorien.Range("J21").Value = "Popkultur"
$wsKategorien.Range("J21").Style = "Gut"
$wsKategorien.Range("K21").Value = "Guillermo"

$wsKategorien.Range("I22").Value = "Movie Themes"
$wsKategorien.Range("J22").Value = "Popkultur"
$wsKategorien.Range("J22").Style = "Gut"

$wsKategorien.Range("I23").Value = "Recursive Acr"
$wsKategorien.Range("J23").Value = "Technisch"
$wsKategorien.Range("J23").Style = "Schlecht"

# ---------------------------------------------------------------------
# Kategorien sheet: Round 4 (new block, rows 25-31).
# ---------------------------------------------------------------------
$wsKategorien.Range("I25").Value = "Runde 4"

$wsKategorien.Range("I26").Value = "Text 2 Song"
$wsKategorien.Range("J26").Value = "Popkultur"
$wsKategorien.Range("J26").Style = "Gut"
$wsKategorien.Range("K26").Value = "Cipi"

$wsKategorien.Range("I27").Value = "App Icons"
$wsKategorien.Range("J27").Value = "Popkultur"
$wsKategorien.Range("J27").Style = "Gut"
$wsKategorien.Range("K27").Value = "Markus"

$wsKategorien.Range("I28").Value = "Way Back Machine"
$wsKategorien.Range("J28").Value = "Popkultur"
$wsKategorien.Range("J28").Style = "Gut"
$wsKategorien.Range("K28").Value = "Tobi"

$wsKategorien.Range("I29").Value = "Who are they"
$wsKategorien.Range("J29").Value = "virtual7"
$wsKategorien.Range("J29").Style = "Neutral"
$wsKategorien.Range("K29").Value = "Dezsö"

$wsKategorien.Range("I30").Value = "Unix Commands"
$wsKategorien.Range("J30").Value = "Technisch"
$wsKategorien.Range("J30").Style = "Schlecht"

$wsKategorien.Range("I31").Value = "Movie Themes"
$wsKategorien.Range("J31").Value = "Popkultur"
$wsKategorien.Range("J31").Style = "Gut"

# ---------------------------------------------------------------------
# Kategorien sheet: Finale block (rows 33-39).
# ---------------------------------------------------------------------
$wsKategorien.Range("I33").Value = "Finale"

$wsKategorien.Range("I34").Value = "Numbers"
$wsKategorien.Range("J34").Value = "Technisch"
$wsKategorien.Range("J34").Style = "Schlecht"
$wsKategorien.Range("K34").Value = "???"

$wsKategorien.Range("I35").Value = "Text 2 Song"
$wsKategorien.Range("J35").Value = "Popkultur"
$wsKategorien.Range("J35").Style = "Gut"

$wsKategorien.Range("I36").Value = "Movie Themes"
$wsKategorien.Range("J36").Value = "Popkultur"
$wsKategorien.Range("J36").Style = "Gut"

$wsKategorien.Range("I37").Value = "Serien Rollen"
$wsKategorien.Range("J37").Value = "Popkultur"
$wsKategorien.Range("J37").Style = "Gut"

$wsKategorien.Range("I38").Value = "Who am I / Who are they"
$wsKategorien.Range("J38").Value = "virtual7"
$wsKategorien.Range("J38").Style = "Neutral"

$wsKategorien.Range("I39").Value = "Oracle Ports"
$wsKategorien.Range("J39").Value = "Technisch"
$wsKategorien.Range("J39").Style = "Schlecht"

# ---------------------------------------------------------------------
# Selections - keep the Kategorien sheet active (as before) and move its
# cursor down to the newly-entered K18 cell; select the whole used range
# on the Spieler sheet.
# ---------------------------------------------------------------------
$wsSpieler.Range("A1:D17").Select()
$wsKategorien.Activate()
$wsKategorien.Range("K18").Select()
